# Update Name of Algo
# Apply updated KNN-imputed values to column E for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 17.386
    27 = 16.472
    32 = 17.4
    36 = 16.652
    38 = 16.631
    46 = 16.966
    54 = 16.754
    55 = 16.354
    56 = 16.351
    67 = 17.24
    69 = 17.386
    72 = 16.939
    83 = 16.886
    86 = 16.542
    91 = 17.577
    93 = 17.22
    99 = 16.659
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
